$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 243, shifting existing rows 243:362 down to 244:363
$ws.Rows("243:243").Insert()

# Populate the newly inserted row 243 with its data.
# Columns A, B, C, E, F, G, H, I, R are identical to the surrounding "Betarraga" rows.
$ws.Range("A243").Value = 4
$ws.Range("B243").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C243").Value = "Los Lagos"
$ws.Range("D243").Value = 44839
$ws.Range("E243").Value = 10
$ws.Range("F243").Value = 100114014
$ws.Range("G243").Value = "Betarraga"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 100
$ws.Range("K243").Value = 1500
$ws.Range("L243").Value = 1500
$ws.Range("M243").Value = 1500
$ws.Range("N243").Value = "$/paquete 5 unidades"
$ws.Range("O243").Value = "Región del Maule"
$ws.Range("P243").Value = 300
$ws.Range("Q243").Value = 5
$ws.Range("R243").Value = "Hortaliza"

# Make sure the new D243 date cell keeps the same date style used by the rest of column D
$ws.Range("D243").NumberFormat = $ws.Range("D244").NumberFormat
